# Time recording log: start a new week ("Nädal 6") by duplicating the
# latest week's sheet ("Nädal 5"), placing the copy in front of it, then
# clearing the copy's logged rows (new week hasn't been filled in yet,
# except for the very first entry) while finishing off the old week's
# last entry (row 16) on the original sheet.

$wb = $excel.ActiveWorkbook

# The most-recently-used week sheet is the first tab.
$latest = $wb.Worksheets.Item(1)
$latestName = $latest.Name

# Duplicate it, inserting the copy before the original -> becomes the new
# first tab / active sheet.
$latest.Copy($latest)

# Fetch stable references by name (worksheet indices shift after Copy).
$newWeek = $wb.Worksheets.Item(1)
$newWeek.Name = "Nädal 6"
$oldWeek = $wb.Worksheets.Item($latestName)

# --- "Nädal 5" sheet (the original data): fill in the final logged row ---
$oldWeek.Range("F15").Value = 25
$oldWeek.Range("D16").Value = 0.95138888888888884
$oldWeek.Range("F16").Value = 370

[void]$oldWeek.Range("F16").Select()

# --- New "Nädal 6" sheet: clear out the previous week's log entries ---
# Row 7 keeps a single carried-over entry, but with fresh date/start time,
# no stop time / delta minutes yet, and a new activity comment.
$newWeek.Range("B7").Value = 43892
$newWeek.Range("C7").Value = 0.36458333333333331
$newWeek.Range("D7").ClearContents()
$newWeek.Range("F7").ClearContents()
$newWeek.Range("G7").Value = "Eelmise kodutöö errorite eelimineerimine"

# Rows 8-16 are fully blanked out - nothing logged yet this week.
$newWeek.Range("B8:H16").ClearContents()

[void]$newWeek.Activate()
[void]$newWeek.Range("G7").Select()
